$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1780343333333333
$ws.Range("H2").Value = 0.534103
$ws.Range("I2").Value = 0.003649670474736916
$ws.Range("J2").Value = 0.003649670474736915
$ws.Range("M2").Value = 0.08962966666666666
$ws.Range("N2").Value = 0.268889
$ws.Range("O2").Value = 0.4339761198462219
$ws.Range("P2").Value = 0.4339761198462219
$ws.Range("Q2").Value = 0.01595715795188889
$ws.Range("R2").Value = 0.143614421567
$ws.Range("S2").Value = 0.001583869831343645
$ws.Range("T2").Value = 0.001583869831343645

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1780343333333333
$ws.Range("H3").Value = 0.534103
$ws.Range("I3").Value = 0.003649670474736916
$ws.Range("J3").Value = 0.003649670474736915
$ws.Range("M3").Value = 0.1169016666666667
$ws.Range("N3").Value = 0.350705
$ws.Range("O3").Value = 0.5660238801537781
$ws.Range("P3").Value = 0.5660238801537781
$ws.Range("Q3").Value = 0.02081251029055555
$ws.Range("R3").Value = 0.187312592615
$ws.Range("S3").Value = 0.002065800643393271
$ws.Range("T3").Value = 0.002065800643393271

# Row 4
$ws.Range("I4").Value = 0.09908483984804967
$ws.Range("J4").Value = 0.09908483984804965
$ws.Range("M4").Value = 0.08962966666666666
$ws.Range("N4").Value = 0.268889
$ws.Range("O4").Value = 0.4339761198462219
$ws.Range("P4").Value = 0.4339761198462219
$ws.Range("Q4").Value = 0.4332206019796667
$ws.Range("R4").Value = 3.898985417817
$ws.Range("S4").Value = 0.0430004543328409
$ws.Range("T4").Value = 0.04300045433284089

# Row 5
$ws.Range("I5").Value = 0.09908483984804967
$ws.Range("J5").Value = 0.09908483984804965
$ws.Range("M5").Value = 0.1169016666666667
$ws.Range("N5").Value = 0.350705
$ws.Range("O5").Value = 0.5660238801537781
$ws.Range("P5").Value = 0.5660238801537781
$ws.Range("Q5").Value = 0.5650384776516667
$ws.Range("R5").Value = 5.085346298865
$ws.Range("S5").Value = 0.05608438551520877
$ws.Range("T5").Value = 0.05608438551520876

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 35.63223
$ws.Range("H6").Value = 106.89669
$ws.Range("I6").Value = 0.7304540385283456
$ws.Range("J6").Value = 0.7304540385283456
$ws.Range("M6").Value = 0.08962966666666666
$ws.Range("N6").Value = 0.268889
$ws.Range("O6").Value = 0.4339761198462219
$ws.Range("P6").Value = 0.4339761198462219
$ws.Range("Q6").Value = 3.19370489749
$ws.Range("R6").Value = 28.74334407741
$ws.Range("S6").Value = 0.316999609366534
$ws.Range("T6").Value = 0.316999609366534

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 35.63223
$ws.Range("H7").Value = 106.89669
$ws.Range("I7").Value = 0.7304540385283456
$ws.Range("J7").Value = 0.7304540385283456
$ws.Range("M7").Value = 0.1169016666666667
$ws.Range("N7").Value = 0.350705
$ws.Range("O7").Value = 0.5660238801537781
$ws.Range("P7").Value = 0.5660238801537781
$ws.Range("Q7").Value = 4.16546707405
$ws.Range("R7").Value = 37.48920366645
$ws.Range("S7").Value = 0.4134544291618115
$ws.Range("T7").Value = 0.4134544291618115

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5521946666666667
$ws.Range("H8").Value = 1.656584
$ws.Range("I8").Value = 0.011319887201011
$ws.Range("J8").Value = 0.011319887201011
$ws.Range("M8").Value = 0.08962966666666666
$ws.Range("N8").Value = 0.268889
$ws.Range("O8").Value = 0.4339761198462219
$ws.Range("P8").Value = 0.4339761198462219
$ws.Range("Q8").Value = 0.04949302390844445
$ws.Range("R8").Value = 0.445437215176
$ws.Range("S8").Value = 0.004912560724591664
$ws.Range("T8").Value = 0.004912560724591663

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5521946666666667
$ws.Range("H9").Value = 1.656584
$ws.Range("I9").Value = 0.011319887201011
$ws.Range("J9").Value = 0.011319887201011
$ws.Range("M9").Value = 0.1169016666666667
$ws.Range("N9").Value = 0.350705
$ws.Range("O9").Value = 0.5660238801537781
$ws.Range("P9").Value = 0.5660238801537781
$ws.Range("Q9").Value = 0.06455247685777779
$ws.Range("R9").Value = 0.58097229172
$ws.Range("S9").Value = 0.006407326476419339
$ws.Range("T9").Value = 0.006407326476419339

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.585023666666667
$ws.Range("H10").Value = 22.755071
$ws.Range("I10").Value = 0.155491563947857
$ws.Range("J10").Value = 0.1554915639478569
$ws.Range("M10").Value = 0.08962966666666666
$ws.Range("N10").Value = 0.268889
$ws.Range("O10").Value = 0.4339761198462219
$ws.Range("P10").Value = 0.4339761198462219
$ws.Range("Q10").Value = 0.6798431429021111
$ws.Range("R10").Value = 6.118588286119
$ws.Range("S10").Value = 0.06747962559091164
$ws.Range("T10").Value = 0.06747962559091163

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.585023666666667
$ws.Range("H11").Value = 22.755071
$ws.Range("I11").Value = 0.155491563947857
$ws.Range("J11").Value = 0.1554915639478569
$ws.Range("M11").Value = 0.1169016666666667
$ws.Range("N11").Value = 0.350705
$ws.Range("O11").Value = 0.5660238801537781
$ws.Range("P11").Value = 0.5660238801537781
$ws.Range("Q11").Value = 0.8867019083394445
$ws.Range("R11").Value = 7.980317175055
$ws.Range("S11").Value = 0.08801193835694532
$ws.Range("T11").Value = 0.0880119383569453
